$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 421.6
$ws.Range("I11").Value = 421.6
$ws.Range("K11").Value = 421.6
$ws.Range("M11").Value = -281.6
$ws.Range("H40").Value = 1194.1177
$ws.Range("I40").Value = 1194.1177
$ws.Range("K40").Value = 1194.1177
$ws.Range("M40").Value = -1019.1177
$ws.Range("H55").Value = 103.44444
$ws.Range("I55").Value = 106.333336
$ws.Range("J55").Value = 102
$ws.Range("K55").Value = 106.333336
$ws.Range("L55").Value = 102
$ws.Range("M55").Value = 107.666664
$ws.Range("N55").Value = -530
$ws.Range("H64").Value = 3012.8235
$ws.Range("I64").Value = 2898.5
$ws.Range("J64").Value = 3287.2
$ws.Range("K64").Value = 2898.5
$ws.Range("L64").Value = 3287.2
$ws.Range("M64").Value = -2650.5
$ws.Range("N64").Value = -3783.2
$ws.Range("H67").Value = 3012.8235
$ws.Range("I67").Value = 2898.5
$ws.Range("J67").Value = 3287.2
$ws.Range("K67").Value = 2898.5
$ws.Range("L67").Value = 3287.2
$ws.Range("M67").Value = -2040.5
$ws.Range("N67").Value = -5003.2
$ws.Range("H106").Value = 2904.0476
$ws.Range("I106").Value = 2749
$ws.Range("J106").Value = 2952.5
$ws.Range("K106").Value = 2749
$ws.Range("L106").Value = 2952.5
$ws.Range("M106").Value = -2118
$ws.Range("N106").Value = -4214.5
$ws.Range("H125").Value = 1040.3334
$ws.Range("I125").Value = 843.8333
$ws.Range("K125").Value = 7594.4997
$ws.Range("M125").Value = -5134.4997
$ws.Range("H141").Value = 2280.805
$ws.Range("I141").Value = 1966.8529
$ws.Range("J141").Value = 3805.7144
$ws.Range("K141").Value = 5900.5587
$ws.Range("L141").Value = 11417.1432
$ws.Range("M141").Value = -720.5587000000005
$ws.Range("N141").Value = -21777.1432

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4855.2
$ws.Range("I32").Value = 3762.0964
$ws.Range("J32").Value = 17816.285
$ws.Range("K32").Value = 3762.0964
$ws.Range("L32").Value = 17816.285
$ws.Range("M32").Value = -3475.0964
$ws.Range("N32").Value = -18390.285
$ws.Range("H102").Value = 3527.5
$ws.Range("I102").Value = 3527.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3527.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1905.5
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 1925.18
$ws.Range("I132").Value = 1911.1945
$ws.Range("J132").Value = 1961.1428
$ws.Range("K132").Value = 5733.583500000001
$ws.Range("L132").Value = 5883.428400000001
$ws.Range("M132").Value = -3203.583500000001
$ws.Range("N132").Value = -10943.4284

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 44524.3
$ws.Range("I31").Value = 81855.84
$ws.Range("J31").Value = 15976.647
$ws.Range("K31").Value = 81855.84
$ws.Range("L31").Value = 15976.647
$ws.Range("M31").Value = -81560.84
$ws.Range("N31").Value = -16566.647
$ws.Range("H34").Value = 44524.3
$ws.Range("I34").Value = 81855.84
$ws.Range("J34").Value = 15976.647
$ws.Range("K34").Value = 81855.84
$ws.Range("L34").Value = 15976.647
$ws.Range("M34").Value = -81653.84
$ws.Range("N34").Value = -16380.647
$ws.Range("H94").Value = 191487.28
$ws.Range("I94").Value = 222870.44
$ws.Range("J94").Value = 167949.92
$ws.Range("K94").Value = 222870.44
$ws.Range("L94").Value = 167949.92
$ws.Range("M94").Value = -222419.44
$ws.Range("N94").Value = -168851.92
$ws.Range("H99").Value = 1583.3334
$ws.Range("I99").Value = 1375
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1375
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 123
$ws.Range("N99").Value = -4996
$ws.Range("H122").Value = 3125.3333
$ws.Range("I122").Value = 3275.3845
$ws.Range("J122").Value = 2150
$ws.Range("K122").Value = 9826.1535
$ws.Range("L122").Value = 6450
$ws.Range("M122").Value = -7376.1535
$ws.Range("N122").Value = -11350
$ws.Range("H126").Value = 1583.3334
$ws.Range("I126").Value = 1375
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 4125
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1655
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 1071.7858
$ws.Range("I132").Value = 855.2857
$ws.Range("J132").Value = 2154.2856
$ws.Range("K132").Value = 2565.8571
$ws.Range("L132").Value = 6462.8568
$ws.Range("M132").Value = -35.85710000000017
$ws.Range("N132").Value = -11522.8568
$ws.Range("H134").Value = 1756.4131
$ws.Range("I134").Value = 1409.2927
$ws.Range("J134").Value = 4602.8
$ws.Range("K134").Value = 4227.8781
$ws.Range("L134").Value = 13808.4
$ws.Range("M134").Value = -1692.8781
$ws.Range("N134").Value = -18878.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4893.5654
$ws.Range("I68").Value = 354.57144
$ws.Range("J68").Value = 6879.375
$ws.Range("K68").Value = 1063.71432
$ws.Range("L68").Value = 20638.125
$ws.Range("M68").Value = -252.71432
$ws.Range("N68").Value = -22260.125
$ws.Range("H71").Value = 4893.5654
$ws.Range("I71").Value = 354.57144
$ws.Range("J71").Value = 6879.375
$ws.Range("K71").Value = 3191.14296
$ws.Range("L71").Value = 61914.375
$ws.Range("M71").Value = 864.8570399999999
$ws.Range("N71").Value = -70026.375
$ws.Range("H107").Value = 8817.25
$ws.Range("I107").Value = 33551
$ws.Range("J107").Value = 572.6667
$ws.Range("K107").Value = 100653
$ws.Range("L107").Value = 1718.0001
$ws.Range("M107").Value = -98733
$ws.Range("N107").Value = -5558.0001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1487.625
$ws.Range("I16").Value = 1025.1666
$ws.Range("J16").Value = 2875
$ws.Range("K16").Value = 1025.1666
$ws.Range("L16").Value = 2875
$ws.Range("M16").Value = -855.1666
$ws.Range("N16").Value = -3215
$ws.Range("H132").Value = 1642.3522
$ws.Range("I132").Value = 1454.7759
$ws.Range("J132").Value = 2479.2307
$ws.Range("K132").Value = 4364.3277
$ws.Range("L132").Value = 7437.6921
$ws.Range("M132").Value = -1834.3277
$ws.Range("N132").Value = -12497.6921
$ws.Range("H136").Value = 1663.134
$ws.Range("I136").Value = 1130.973
$ws.Range("K136").Value = 3392.919
$ws.Range("M136").Value = -842.9189999999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 379.69766
$ws.Range("I136").Value = 294.13333
$ws.Range("J136").Value = 577.1539
$ws.Range("K136").Value = 882.39999
$ws.Range("L136").Value = 1731.4617
$ws.Range("M136").Value = 1667.60001
$ws.Range("N136").Value = -6831.4617
